# "Running all the test cases in Chrome"
# Test cases B2-B6 (rows 3 through 7 on the "Test Cases" sheet) were re-run
# and now pass ("Y") instead of the previous fail result ("N").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$ws.Range("C3:C7").Value = "Y"
